$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: rowNumber, Coin, Link, Price, Volume(1h), priceIsNumericLooking, volIsNumericLooking
# "*IsNumericLooking" rows get a leading apostrophe so Excel keeps them as text
# instead of silently parsing them into a double (which would mangle values like
# "0.9960" -> 0.996 or "1.002" -> 1.002 as a real number instead of literal text).
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '20.747.75', '  +2.40%  ', 0, 0)
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.511.34', '  +4.23%  ', 0, 0)
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.002', '  -0.40%  ', 1, 0)
    ,@(5, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9634', '  +2.22%  ', 1, 0)
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '279.44', '  +1.97%  ', 1, 0)
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.3562', '  -1.91%  ', 1, 0)
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3108', '  +1.46%  ', 1, 0)
    ,@(9, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.098', '  +6.67%  ', 1, 0)
    ,@(10, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '39.44', '  -0.58%  ', 1, 0)
    ,@(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06717', '  +3.07%  ', 1, 0)
    ,@(12, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9960', '  -0.18%  ', 1, 0)
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '18.38', '  +4.49%  ', 1, 0)
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.533', '  +3.36%  ', 1, 0)
    ,@(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.230', '  +2.64%  ', 1, 0)
    ,@(16, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9608', '  +0.40%  ', 1, 0)
    ,@(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001024', '  +1.29%  ', 1, 0)
    ,@(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.484.33', '  +2.71%  ', 0, 0)
    ,@(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06019', '  +5.72%  ', 1, 0)
    ,@(20, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '69.81', '  +1.37%  ', 1, 0)
    ,@(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.547', '  +2.94%  ', 1, 0)
    ,@(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.65', '  +2.29%  ', 1, 0)
    ,@(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.19', '  +3.78%  ', 1, 0)
    ,@(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.303', '  +3.21%  ', 1, 0)
    ,@(25, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '20.692.39', '  +1.98%  ', 0, 0)
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '146.18', '  +4.01%  ', 1, 0)
    ,@(27, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.126', '  +1.99%  ', 1, 0)
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.43', '  +2.79%  ', 1, 0)
    ,@(29, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.636.46', '  +2.48%  ', 0, 0)
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '115.67', '  +4.21%  ', 1, 0)
    ,@(31, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.977', '  +0.62%  ', 1, 0)
    ,@(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.032', '  +4.62%  ', 1, 0)
    ,@(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.8141', '  +3.51%  ', 1, 0)
    ,@(34, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.07965', '  +3.28%  ', 1, 0)
    ,@(35, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.203', '  +7.96%  ', 1, 0)
    ,@(36, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.435', '  -2.91%  ', 1, 0)
    ,@(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05757', '  +2.02%  ', 1, 0)
    ,@(38, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.793', '  +2.95%  ', 1, 0)
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02040', '  +1.97%  ', 1, 0)
    ,@(40, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.44', '  +2.36%  ', 1, 0)
    ,@(41, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9617', '  +1.41%  ', 1, 0)
    ,@(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '7.605', '  +3.14%  ', 1, 0)
    ,@(43, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1864', '  +0.64%  ', 1, 0)
    ,@(44, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5283', '  +1.20%  ', 1, 0)
    ,@(45, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.534', '  +1.75%  ', 1, 0)
    ,@(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.04', '  +1.17%  ', 1, 0)
    ,@(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '119.93', '  +3.18%  ', 1, 0)
    ,@(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5234', '  +2.55%  ', 1, 0)
    ,@(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.834', '  +5.91%  ', 1, 0)
    ,@(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06466', '  +1.35%  ', 1, 0)
    ,@(51, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9896', '  +0.68%  ', 1, 0)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value2 = $r[1]
    $ws.Cells.Item($row, 3).Value2 = $r[2]

    if ($r[5] -eq 1) {
        $ws.Cells.Item($row, 4).Value2 = "'" + $r[3]
        $ws.Cells.Item($row, 4).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 4).Value2 = $r[3]
    }

    if ($r[6] -eq 1) {
        $ws.Cells.Item($row, 5).Value2 = "'" + $r[4]
        $ws.Cells.Item($row, 5).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 5).Value2 = $r[4]
    }
}
